$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.890.39"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "2.346.94"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'548.45"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "'131.72"
$ws.Range("E6").Value = "  -0.95%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").Value = "2.348.44"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").Value = "'5.53"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").Value = "'0.338"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "'23.64"
$ws.Range("E14").Value = "  -1.85%  "

$ws.Range("D15").Value = "2.766.17"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").Value = "60.832.18"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").Value = "2.340.44"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").Value = "'10.67"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.10"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'315.56"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'6.63"
$ws.Range("E22").Value = "  -3.85%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'64.29"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("D25").Value = "'0.173"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'7.95"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").Value = "'1.41"
$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("D29").Value = "'1.26"
$ws.Range("E29").Value = "  +9.44%  "

$ws.Range("D30").Value = "'172.40"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("D32").Value = "0.0₃0737"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").Value = "'5.98"
$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("D34").Value = "'1.38"
$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("D35").Value = "'0.385"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "'18.02"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").Value = "'4.18"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "'326.82"
$ws.Range("E40").Value = "  +3.10%  "

$ws.Range("D41").Value = "'38.26"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "'1.54"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").Value = "'137.48"
$ws.Range("E43").Value = "  -4.10%  "

$ws.Range("D44").Value = "'3.51"
$ws.Range("E44").Value = "  +0.71%  "

$ws.Range("D45").Value = "'0.0945"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").Value = "'19.31"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("D47").Value = "'0.573"
$ws.Range("E47").Value = "  +1.54%  "

$ws.Range("D48").Value = "'0.0498"
$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("D49").Value = "'0.0217"
$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("D50").Value = "0.0⁦0221"
$ws.Range("E50").Value = "  +4.30%  "

$ws.Range("D51").Value = "'11.00"
$ws.Range("E51").Value = "  -0.35%  "
